# Natmi following Dr Hou advice
# Updates the Vegfb-Flt1 ligand-receptor communication table: the per-cluster
# cell counts were recomputed (Ligand/Receptor-expressing cell counts and all
# of the derived expression / specificity statistics that depend on them),
# and the sheet grows from 12 data rows (3 target clusters) to 16 data rows
# (4 target clusters: ECs, FAPs, M2, sCs) for every sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Cells.Item(2, 5).Value = 3
    $ws.Cells.Item(2, 7).Value = 1.80847
    $ws.Cells.Item(2, 8).Value = 5.42541
    $ws.Cells.Item(2, 9).Value = 0.1042562806587584
    $ws.Cells.Item(2, 10).Value = 0.1042562806587584
    $ws.Cells.Item(2, 11).Value = 3
    $ws.Cells.Item(2, 13).Value = 100.95625
    $ws.Cells.Item(2, 14).Value = 302.86875
    $ws.Cells.Item(2, 15).Value = 0.9361599519103526
    $ws.Cells.Item(2, 16).Value = 0.9361599519103525
    $ws.Cells.Item(2, 17).Value = 182.5763494375
    $ws.Cells.Item(2, 18).Value = 1643.1871449375
    $ws.Cells.Item(2, 19).Value = 0.09760055468785553
    $ws.Cells.Item(2, 20).Value = 0.09760055468785551
    $ws.Cells.Item(3, 4).Value = "FAPs"
    $ws.Cells.Item(3, 5).Value = 3
    $ws.Cells.Item(3, 7).Value = 1.80847
    $ws.Cells.Item(3, 8).Value = 5.42541
    $ws.Cells.Item(3, 9).Value = 0.1042562806587584
    $ws.Cells.Item(3, 10).Value = 0.1042562806587584
    $ws.Cells.Item(3, 11).Value = 2
    $ws.Cells.Item(3, 12).Value = 0.6666666666666666
    $ws.Cells.Item(3, 13).Value = 0.3890956666666667
    $ws.Cells.Item(3, 14).Value = 1.167287
    $ws.Cells.Item(3, 15).Value = 0.003608055772626195
    $ws.Cells.Item(3, 16).Value = 0.003608055772626195
    $ws.Cells.Item(3, 17).Value = 0.7036678402966667
    $ws.Cells.Item(3, 18).Value = 6.33301056267
    $ws.Cells.Item(3, 19).Value = 0.0003761624752633701
    $ws.Cells.Item(3, 20).Value = 0.0003761624752633701
    $ws.Cells.Item(4, 4).Value = "M2"
    $ws.Cells.Item(4, 5).Value = 3
    $ws.Cells.Item(4, 7).Value = 1.80847
    $ws.Cells.Item(4, 8).Value = 5.42541
    $ws.Cells.Item(4, 9).Value = 0.1042562806587584
    $ws.Cells.Item(4, 10).Value = 0.1042562806587584
    $ws.Cells.Item(4, 11).Value = 3
    $ws.Cells.Item(4, 13).Value = 2.996608333333334
    $ws.Cells.Item(4, 14).Value = 8.989825000000002
    $ws.Cells.Item(4, 15).Value = 0.02778733078167519
    $ws.Cells.Item(4, 16).Value = 0.02778733078167519
    $ws.Cells.Item(4, 17).Value = 5.419276272583335
    $ws.Cells.Item(4, 18).Value = 48.77348645325001
    $ws.Cells.Item(4, 19).Value = 0.002897003756732086
    $ws.Cells.Item(4, 20).Value = 0.002897003756732086
    $ws.Cells.Item(5, 1).Value = "ECs"
    $ws.Cells.Item(5, 2).Value = "Vegfb"
    $ws.Cells.Item(5, 4).Value = "sCs"
    $ws.Cells.Item(5, 5).Value = 3
    $ws.Cells.Item(5, 7).Value = 1.80847
    $ws.Cells.Item(5, 8).Value = 5.42541
    $ws.Cells.Item(5, 9).Value = 0.1042562806587584
    $ws.Cells.Item(5, 10).Value = 0.1042562806587584
    $ws.Cells.Item(5, 11).Value = 3
    $ws.Cells.Item(5, 13).Value = 3.498858666666667
    $ws.Cells.Item(5, 14).Value = 10.496576
    $ws.Cells.Item(5, 15).Value = 0.03244466153534613
    $ws.Cells.Item(5, 16).Value = 0.03244466153534613
    $ws.Cells.Item(5, 17).Value = 6.327580932906668
    $ws.Cells.Item(5, 18).Value = 56.94822839616001
    $ws.Cells.Item(5, 19).Value = 0.003382559738907471
    $ws.Cells.Item(5, 20).Value = 0.003382559738907471
    $ws.Cells.Item(6, 4).Value = "ECs"
    $ws.Cells.Item(6, 5).Value = 3
    $ws.Cells.Item(6, 7).Value = 6.255752333333334
    $ws.Cells.Item(6, 8).Value = 18.767257
    $ws.Cells.Item(6, 9).Value = 0.360637152397155
    $ws.Cells.Item(6, 10).Value = 0.360637152397155
    $ws.Cells.Item(6, 11).Value = 3
    $ws.Cells.Item(6, 13).Value = 100.95625
    $ws.Cells.Item(6, 14).Value = 302.86875
    $ws.Cells.Item(6, 15).Value = 0.9361599519103526
    $ws.Cells.Item(6, 16).Value = 0.9361599519103525
    $ws.Cells.Item(6, 17).Value = 631.5572965020833
    $ws.Cells.Item(6, 18).Value = 5684.015668518749
    $ws.Cells.Item(6, 19).Value = 0.3376140592452072
    $ws.Cells.Item(6, 20).Value = 0.3376140592452072
    $ws.Cells.Item(7, 4).Value = "FAPs"
    $ws.Cells.Item(7, 5).Value = 3
    $ws.Cells.Item(7, 7).Value = 6.255752333333334
    $ws.Cells.Item(7, 8).Value = 18.767257
    $ws.Cells.Item(7, 9).Value = 0.360637152397155
    $ws.Cells.Item(7, 10).Value = 0.360637152397155
    $ws.Cells.Item(7, 11).Value = 2
    $ws.Cells.Item(7, 12).Value = 0.6666666666666666
    $ws.Cells.Item(7, 13).Value = 0.3890956666666667
    $ws.Cells.Item(7, 14).Value = 1.167287
    $ws.Cells.Item(7, 15).Value = 0.003608055772626195
    $ws.Cells.Item(7, 16).Value = 0.003608055772626195
    $ws.Cells.Item(7, 17).Value = 2.434086124639889
    $ws.Cells.Item(7, 18).Value = 21.906775121759
    $ws.Cells.Item(7, 19).Value = 0.001301198959530028
    $ws.Cells.Item(7, 20).Value = 0.001301198959530028
    $ws.Cells.Item(8, 1).Value = "FAPs"
    $ws.Cells.Item(8, 2).Value = "Vegfb"
    $ws.Cells.Item(8, 4).Value = "M2"
    $ws.Cells.Item(8, 5).Value = 3
    $ws.Cells.Item(8, 7).Value = 6.255752333333334
    $ws.Cells.Item(8, 8).Value = 18.767257
    $ws.Cells.Item(8, 9).Value = 0.360637152397155
    $ws.Cells.Item(8, 10).Value = 0.360637152397155
    $ws.Cells.Item(8, 11).Value = 3
    $ws.Cells.Item(8, 13).Value = 2.996608333333334
    $ws.Cells.Item(8, 14).Value = 8.989825000000002
    $ws.Cells.Item(8, 15).Value = 0.02778733078167519
    $ws.Cells.Item(8, 16).Value = 0.02778733078167519
    $ws.Cells.Item(8, 17).Value = 18.74603957333612
    $ws.Cells.Item(8, 18).Value = 168.714356160025
    $ws.Cells.Item(8, 19).Value = 0.01002114384582115
    $ws.Cells.Item(8, 20).Value = 0.01002114384582115
    $ws.Cells.Item(9, 1).Value = "FAPs"
    $ws.Cells.Item(9, 2).Value = "Vegfb"
    $ws.Cells.Item(9, 4).Value = "sCs"
    $ws.Cells.Item(9, 5).Value = 3
    $ws.Cells.Item(9, 7).Value = 6.255752333333334
    $ws.Cells.Item(9, 8).Value = 18.767257
    $ws.Cells.Item(9, 9).Value = 0.360637152397155
    $ws.Cells.Item(9, 10).Value = 0.360637152397155
    $ws.Cells.Item(9, 11).Value = 3
    $ws.Cells.Item(9, 13).Value = 3.498858666666667
    $ws.Cells.Item(9, 14).Value = 10.496576
    $ws.Cells.Item(9, 15).Value = 0.03244466153534613
    $ws.Cells.Item(9, 16).Value = 0.03244466153534613
    $ws.Cells.Item(9, 17).Value = 21.88799326800356
    $ws.Cells.Item(9, 18).Value = 196.991939412032
    $ws.Cells.Item(9, 19).Value = 0.01170075034659674
    $ws.Cells.Item(9, 20).Value = 0.01170075034659674
    $ws.Cells.Item(10, 4).Value = "ECs"
    $ws.Cells.Item(10, 5).Value = 3
    $ws.Cells.Item(10, 7).Value = 3.749018666666667
    $ws.Cells.Item(10, 8).Value = 11.247056
    $ws.Cells.Item(10, 9).Value = 0.2161267493001954
    $ws.Cells.Item(10, 10).Value = 0.2161267493001954
    $ws.Cells.Item(10, 11).Value = 3
    $ws.Cells.Item(10, 13).Value = 100.95625
    $ws.Cells.Item(10, 14).Value = 302.86875
    $ws.Cells.Item(10, 15).Value = 0.9361599519103526
    $ws.Cells.Item(10, 16).Value = 0.9361599519103525
    $ws.Cells.Item(10, 17).Value = 378.4868657666667
    $ws.Cells.Item(10, 18).Value = 3406.3817919
    $ws.Cells.Item(10, 19).Value = 0.2023292072314118
    $ws.Cells.Item(10, 20).Value = 0.2023292072314117
    $ws.Cells.Item(11, 1).Value = "M2"
    $ws.Cells.Item(11, 2).Value = "Vegfb"
    $ws.Cells.Item(11, 4).Value = "FAPs"
    $ws.Cells.Item(11, 5).Value = 3
    $ws.Cells.Item(11, 7).Value = 3.749018666666667
    $ws.Cells.Item(11, 8).Value = 11.247056
    $ws.Cells.Item(11, 9).Value = 0.2161267493001954
    $ws.Cells.Item(11, 10).Value = 0.2161267493001954
    $ws.Cells.Item(11, 11).Value = 2
    $ws.Cells.Item(11, 12).Value = 0.6666666666666666
    $ws.Cells.Item(11, 13).Value = 0.3890956666666667
    $ws.Cells.Item(11, 14).Value = 1.167287
    $ws.Cells.Item(11, 15).Value = 0.003608055772626195
    $ws.Cells.Item(11, 16).Value = 0.003608055772626195
    $ws.Cells.Item(11, 17).Value = 1.458726917452444
    $ws.Cells.Item(11, 18).Value = 13.128542257072
    $ws.Cells.Item(11, 19).Value = 0.0007797973654315044
    $ws.Cells.Item(11, 20).Value = 0.0007797973654315043
    $ws.Cells.Item(12, 1).Value = "M2"
    $ws.Cells.Item(12, 2).Value = "Vegfb"
    $ws.Cells.Item(12, 5).Value = 3
    $ws.Cells.Item(12, 7).Value = 3.749018666666667
    $ws.Cells.Item(12, 8).Value = 11.247056
    $ws.Cells.Item(12, 9).Value = 0.2161267493001954
    $ws.Cells.Item(12, 10).Value = 0.2161267493001954
    $ws.Cells.Item(12, 11).Value = 3
    $ws.Cells.Item(12, 13).Value = 2.996608333333334
    $ws.Cells.Item(12, 14).Value = 8.989825000000002
    $ws.Cells.Item(12, 15).Value = 0.02778733078167519
    $ws.Cells.Item(12, 16).Value = 0.02778733078167519
    $ws.Cells.Item(12, 17).Value = 11.23434057835556
    $ws.Cells.Item(12, 18).Value = 101.1090652052
    $ws.Cells.Item(12, 19).Value = 0.006005585473572717
    $ws.Cells.Item(12, 20).Value = 0.006005585473572716
    $ws.Cells.Item(13, 1).Value = "M2"
    $ws.Cells.Item(13, 2).Value = "Vegfb"
    $ws.Cells.Item(13, 5).Value = 3
    $ws.Cells.Item(13, 7).Value = 3.749018666666667
    $ws.Cells.Item(13, 8).Value = 11.247056
    $ws.Cells.Item(13, 9).Value = 0.2161267493001954
    $ws.Cells.Item(13, 10).Value = 0.2161267493001954
    $ws.Cells.Item(13, 11).Value = 3
    $ws.Cells.Item(13, 13).Value = 3.498858666666667
    $ws.Cells.Item(13, 14).Value = 10.496576
    $ws.Cells.Item(13, 15).Value = 0.03244466153534613
    $ws.Cells.Item(13, 16).Value = 0.03244466153534613
    $ws.Cells.Item(13, 17).Value = 13.11728645336178
    $ws.Cells.Item(13, 18).Value = 118.055578080256
    $ws.Cells.Item(13, 19).Value = 0.007012159229779446
    $ws.Cells.Item(13, 20).Value = 0.007012159229779445
    $ws.Cells.Item(14, 1).Value = "sCs"
    $ws.Cells.Item(14, 2).Value = "Vegfb"
    $ws.Cells.Item(14, 3).Value = "Flt1"
    $ws.Cells.Item(14, 4).Value = "ECs"
    $ws.Cells.Item(14, 5).Value = 3
    $ws.Cells.Item(14, 6).Value = 1
    $ws.Cells.Item(14, 7).Value = 5.533148
    $ws.Cells.Item(14, 8).Value = 16.599444
    $ws.Cells.Item(14, 9).Value = 0.3189798176438912
    $ws.Cells.Item(14, 10).Value = 0.3189798176438912
    $ws.Cells.Item(14, 11).Value = 3
    $ws.Cells.Item(14, 12).Value = 1
    $ws.Cells.Item(14, 13).Value = 100.95625
    $ws.Cells.Item(14, 14).Value = 302.86875
    $ws.Cells.Item(14, 15).Value = 0.9361599519103526
    $ws.Cells.Item(14, 16).Value = 0.9361599519103525
    $ws.Cells.Item(14, 17).Value = 558.605872775
    $ws.Cells.Item(14, 18).Value = 5027.452854974999
    $ws.Cells.Item(14, 19).Value = 0.2986161307458782
    $ws.Cells.Item(14, 20).Value = 0.2986161307458782
    $ws.Cells.Item(15, 1).Value = "sCs"
    $ws.Cells.Item(15, 2).Value = "Vegfb"
    $ws.Cells.Item(15, 3).Value = "Flt1"
    $ws.Cells.Item(15, 4).Value = "FAPs"
    $ws.Cells.Item(15, 5).Value = 3
    $ws.Cells.Item(15, 6).Value = 1
    $ws.Cells.Item(15, 7).Value = 5.533148
    $ws.Cells.Item(15, 8).Value = 16.599444
    $ws.Cells.Item(15, 9).Value = 0.3189798176438912
    $ws.Cells.Item(15, 10).Value = 0.3189798176438912
    $ws.Cells.Item(15, 11).Value = 2
    $ws.Cells.Item(15, 12).Value = 0.6666666666666666
    $ws.Cells.Item(15, 13).Value = 0.3890956666666667
    $ws.Cells.Item(15, 14).Value = 1.167287
    $ws.Cells.Item(15, 15).Value = 0.003608055772626195
    $ws.Cells.Item(15, 16).Value = 0.003608055772626195
    $ws.Cells.Item(15, 17).Value = 2.152923909825333
    $ws.Cells.Item(15, 18).Value = 19.376315188428
    $ws.Cells.Item(15, 19).Value = 0.001150896972401293
    $ws.Cells.Item(15, 20).Value = 0.001150896972401292
    $ws.Cells.Item(16, 1).Value = "sCs"
    $ws.Cells.Item(16, 2).Value = "Vegfb"
    $ws.Cells.Item(16, 3).Value = "Flt1"
    $ws.Cells.Item(16, 4).Value = "M2"
    $ws.Cells.Item(16, 5).Value = 3
    $ws.Cells.Item(16, 6).Value = 1
    $ws.Cells.Item(16, 7).Value = 5.533148
    $ws.Cells.Item(16, 8).Value = 16.599444
    $ws.Cells.Item(16, 9).Value = 0.3189798176438912
    $ws.Cells.Item(16, 10).Value = 0.3189798176438912
    $ws.Cells.Item(16, 11).Value = 3
    $ws.Cells.Item(16, 12).Value = 1
    $ws.Cells.Item(16, 13).Value = 2.996608333333334
    $ws.Cells.Item(16, 14).Value = 8.989825000000002
    $ws.Cells.Item(16, 15).Value = 0.02778733078167519
    $ws.Cells.Item(16, 16).Value = 0.02778733078167519
    $ws.Cells.Item(16, 17).Value = 16.58067740636667
    $ws.Cells.Item(16, 18).Value = 149.2260966573
    $ws.Cells.Item(16, 19).Value = 0.008863597705549236
    $ws.Cells.Item(16, 20).Value = 0.008863597705549236
    $ws.Cells.Item(17, 1).Value = "sCs"
    $ws.Cells.Item(17, 2).Value = "Vegfb"
    $ws.Cells.Item(17, 3).Value = "Flt1"
    $ws.Cells.Item(17, 4).Value = "sCs"
    $ws.Cells.Item(17, 5).Value = 3
    $ws.Cells.Item(17, 6).Value = 1
    $ws.Cells.Item(17, 7).Value = 5.533148
    $ws.Cells.Item(17, 8).Value = 16.599444
    $ws.Cells.Item(17, 9).Value = 0.3189798176438912
    $ws.Cells.Item(17, 10).Value = 0.3189798176438912
    $ws.Cells.Item(17, 11).Value = 3
    $ws.Cells.Item(17, 12).Value = 1
    $ws.Cells.Item(17, 13).Value = 3.498858666666667
    $ws.Cells.Item(17, 14).Value = 10.496576
    $ws.Cells.Item(17, 15).Value = 0.03244466153534613
    $ws.Cells.Item(17, 16).Value = 0.03244466153534613
    $ws.Cells.Item(17, 17).Value = 19.35970283374933
    $ws.Cells.Item(17, 18).Value = 174.237325503744
    $ws.Cells.Item(17, 19).Value = 0.01034919222006248
    $ws.Cells.Item(17, 20).Value = 0.01034919222006248
